# Auto update: 2025-12-01 14:08:26
# Refresh the daily 방산(defense) stock analysis table with the latest
# pulled metrics. Rows 4-6 also change order because HANWHA SYSTEMS now
# ranks ahead of LIG Nex1 and HANWHA AEROSPACE in this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: KOREA AEROSPACE / 047810.KS
$ws.Range("B2").Value = "KOREA AEROSPACE"
$ws.Range("C2").Value = "047810.KS"
$ws.Range("D2").Value = 108300
$ws.Range("E2").Value = 60.8
$ws.Range("F2").Value = 4.23
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 56
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 67.8
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.92500513438651
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# Row 3: HYUNDAI ROTEM / 064350.KS
$ws.Range("B3").Value = "HYUNDAI ROTEM"
$ws.Range("C3").Value = "064350.KS"
$ws.Range("D3").Value = 169900
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = -2.91
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 55.2
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.92500513438651
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"

# Row 4: HANWHA SYSTEMS / 272210.KS (moved up from row 6)
$ws.Range("B4").Value = "HANWHA SYSTEMS"
$ws.Range("C4").Value = "272210.KS"
$ws.Range("D4").Value = 45500
$ws.Range("E4").Value = 19.9
$ws.Range("F4").Value = -4.01
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 66
$ws.Range("J4").Value = 46
$ws.Range("K4").Value = 55.2
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 85.92500513438651
$ws.Range("O4").Value = "🟢 완화적 (상승 우위)"

# Row 5: LIG Nex1 / 079550.KS (moved down from row 4)
$ws.Range("B5").Value = "LIG Nex1"
$ws.Range("C5").Value = "079550.KS"
$ws.Range("D5").Value = 365500
$ws.Range("E5").Value = 33.8
$ws.Range("F5").Value = -7.23
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 56
$ws.Range("K5").Value = 52.8
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 85.92500513438651
$ws.Range("O5").Value = "🟢 완화적 (상승 우위)"

# Row 6: HANWHA AEROSPACE / 012450.KS (moved down from row 5)
$ws.Range("B6").Value = "HANWHA AEROSPACE"
$ws.Range("C6").Value = "012450.KS"
$ws.Range("D6").Value = 812000
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = -5.03
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 60
$ws.Range("J6").Value = 63
$ws.Range("K6").Value = 52.8
$ws.Range("L6").Value = "Pattern"
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 85.92500513438651
$ws.Range("O6").Value = "🟢 완화적 (상승 우위)"

Write-Output "Applied 2025-12-01 refresh to rows 2-6"
